# Fix the straight apostrophes (') to typographic right single quotation
# marks (') in the two CV bullet points that mention "API's" and
# "McDonald's" — matches the author's "Update Docs" text tidy-up.

$d = $word.ActiveDocument
$apos = [char]8217

$d.Content.Find.Execute(
    "Implementing API's from a Java Play framework back-end",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implementing API" + $apos + "s from a Java Play framework back-end",
    2) | Out-Null

$d.Content.Find.Execute(
    "Development of ``Scheduling`` app for McDonald's using Angular and Firebase, to be used in-house for the design team to organise assets scheduled for in-store screen displays",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Development of ``Scheduling`` app for McDonald" + $apos + "s using Angular and Firebase, to be used in-house for the design team to organise assets scheduled for in-store screen displays",
    2) | Out-Null
